# "update at end semester"
# Adds Sprint 4, Sprint 5, Sprint 6 sections (rows 31-57) to the time-tracking
# sheet, mirroring the layout/style of the existing Sprint 1-3 sections, plus
# a semester grand-total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# -----------------------------------------------------------------------
# Phase 1 - write every brand-new piece of text in the exact order it was
# first typed so the shared-string table's append order lines up (existing
# strings like "Coding"/"Testing "/etc. are reused automatically wherever
# else they show up below, regardless of ordering).
# -----------------------------------------------------------------------
$ws.Range("A31").Value = "Sprint 4"
$ws.Range("A40").Value = "Sprint 5"
$ws.Range("A49").Value = "Sprint 6"
$ws.Range("C52").Value = "connected classes between them"
$ws.Range("C42").Value = "Debuggin and reasearching "
$ws.Range("C44").Value = "changed the integration between the classes"
$ws.Range("C51").Value = "changed the design or layer class with a fixed number  of layers (5)"
$ws.Range("A57").Value = "TOTAL SEMESTER  HOURS"
$ws.Range("C35").Value = "Redesigned the component of the layer class"

# -----------------------------------------------------------------------
# Phase 2 - fill in the remaining cell values (all re-use existing shared
# strings), wire up formulas, then restyle/merge everything to match the
# look of the Sprint 1-3 blocks above.
# -----------------------------------------------------------------------

function Stamp-Header($row, $text) {
    # Sprint-header cells (A1/A14/A22) carry a bold/centered style that also
    # spans the two blank cells merged alongside them. Re-using PasteSpecial
    # (instead of poking Font/Alignment directly) makes the engine reuse the
    # existing cellXfs entry instead of minting a near-duplicate one.
    $a = "A" + $row
    $b = "B" + $row
    $c = "C" + $row
    $ws.Range($a).Value = $text
    $ws.Range($b).Value = ""
    $ws.Range($c).Value = ""
    $ws.Range("A1").Copy()
    $ws.Range($a).PasteSpecial($xlPasteFormats)
    $ws.Range($b).PasteSpecial($xlPasteFormats)
    $ws.Range($c).PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false
    $ws.Range($a + ":" + $c).Merge()
}

function Stamp-ColumnHeaders($row) {
    $a = "A" + $row
    $c = "C" + $row
    $ws.Range($a).Value = "Work Type"
    $ws.Range("B" + $row).Value = "Hours Worked"
    $ws.Range($c).Value = "Notes"
    $ws.Range("A2:C2").Copy()
    $ws.Range($a + ":" + $c).PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false
}

function Stamp-DataRow($row, $workType, $hours, $notes, $wrapNotes) {
    $a = "A" + $row
    $b = "B" + $row
    $c = "C" + $row
    $ws.Range($a).Value = $workType
    $ws.Range($b).Value = $hours
    $ws.Range($c).Value = $notes
    $ws.Range("A3:C3").Copy()
    $ws.Range($a + ":" + $c).PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false
    if ($wrapNotes) {
        $ws.Range("C16").Copy()
        $ws.Range($c).PasteSpecial($xlPasteFormats)
        $ws.Application.CutCopyMode = $false
    }
}

function Stamp-TotalRow($row, $label, $formula) {
    $a = "A" + $row
    $b = "B" + $row
    $ws.Range($a).Value = $label
    $ws.Range("A13").Copy()
    $ws.Range($a).PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false
    $ws.Range($b).Formula = $formula
}

# ---------------------------------------------------------------------
# Sprint 4 (rows 31-38)
# ---------------------------------------------------------------------
Stamp-Header 31 "Sprint 4"
Stamp-ColumnHeaders 32
Stamp-DataRow 33 "Coding" 2.5 "Implemented a method to visually hide the the layer from the layer panel " $false
Stamp-DataRow 34 "Coding/Debuging " 4 "fix couple of issues" $false
Stamp-DataRow 35 "Coding/Design" 2 "Redesigned the component of the layer class" $true
Stamp-DataRow 36 "Testing " 1.5 "Debugging and researching " $false
Stamp-TotalRow 38 "Total Hours:" "=B47"

# ---------------------------------------------------------------------
# Sprint 5 (rows 40-47)
# ---------------------------------------------------------------------
Stamp-Header 40 "Sprint 5"
Stamp-ColumnHeaders 41
Stamp-DataRow 42 "Coding" 2 "Debuggin and reasearching " $false
Stamp-DataRow 43 "Coding/Debuging " 1.5 "fix couple of issues" $false
Stamp-DataRow 44 "Coding/Design" 2.5 "changed the integration between the classes" $true
Stamp-DataRow 45 "Testing " 2 "Debugging and researching " $false
Stamp-TotalRow 47 "Total Hours:" "=B42+B43+B44+B45"

# ---------------------------------------------------------------------
# Sprint 6 (rows 49-57) + semester grand total
# ---------------------------------------------------------------------
Stamp-Header 49 "Sprint 6"
Stamp-ColumnHeaders 50
Stamp-DataRow 51 "Coding" 4 "changed the design or layer class with a fixed number  of layers (5)" $false
Stamp-DataRow 52 "Coding/Debuging " 4.5 "connected classes between them" $false
Stamp-DataRow 53 "Coding/Design" 2 "Used Omar's new floating window class to make the layer panel floating and created a docksite for the panel" $true
Stamp-DataRow 54 "Testing " 4 "Debugging and researching " $false
Stamp-TotalRow 56 "Total Hours:" "=B51+B52+B53+B54"
Stamp-TotalRow 57 "TOTAL SEMESTER  HOURS" "=B13+B21+B29+B38+B47+B56"

$ws.Range("C36").Select()

Write-Output "done"
